# Forgot to add dependencies
# Adds the missing second sequencing-run row for sample NA24695 and
# disambiguates the existing output file names with their run ID.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing row 3 (first NA24695 run) output file names now include the
# run identifier to avoid clashing with the new run added below.
$ws.Range("D3").Value = "NA24695_HB66DADXX.g.vcf.gz"
$ws.Range("E3").Value = "NA24695_HB66DADXX.metrics.txt"

# New row 4: second sequencing run (HB7AUADXX) for sample NA24695.
$ws.Range("A4").Value = "ftp://ftp-trace.ncbi.nih.gov/giab/ftp/data/ChineseTrio/HG007_NA24695-hu38168_mother/NA24695_Mother_HiSeq100x/NA24695_Mother_HiSeq100x_fastqs/141117_D00360_0066_BHB7AUADXX/Sample_NA24695/NA24695_CTTGTA_L001_R1_001.fastq.gz"
$ws.Range("B4").Value = "ftp://ftp-trace.ncbi.nih.gov/giab/ftp/data/ChineseTrio/HG007_NA24695-hu38168_mother/NA24695_Mother_HiSeq100x/NA24695_Mother_HiSeq100x_fastqs/141117_D00360_0066_BHB7AUADXX/Sample_NA24695/NA24695_CTTGTA_L001_R2_001.fastq.gz"
$ws.Range("C4").Value = "@RG\tID:HB7AUADXX\tSM:NA24695\tPL:ILLUMINA"
$ws.Range("D4").Value = "NA24695_HB7AUADXX.g.vcf.gz"
$ws.Range("E4").Value = "NA24695_HB7AUADXX.metrics.txt"

# Move the active selection down to reflect the newly appended data.
$ws.Range("E14").Select() | Out-Null
